$d = $word.ActiveDocument

# Locate the paragraph containing the date line "26 - 01 - 2023"
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*26 - 01 - 2023*") {
        $datePara = $p
        break
    }
}

$paraStart = $datePara.Range.Start
$paraText = $datePara.Range.Text

# Find the hyphen that separates "01" from "2023" (the second hyphen).
$dashOffset = $paraText.IndexOf("- 2023")

$dash = $d.Range($paraStart + $dashOffset, $paraStart + $dashOffset + 1)
$dash.Text = "–"

# Force the replaced dash into its own run (distinct from the surrounding
# text runs) by toggling Bold off then back on; this splits the original
# single run into three runs without leaving any stray formatting
# difference behind in the middle run.
$dash.Bold = 0
$dash.Bold = 1

# Add a new right-aligned, bold paragraph right after the date line with
# the professor's name (inherits the date paragraph's formatting).
$datePara.Range.InsertParagraphAfter()
$namePara = $d.Paragraphs.Item($datePara.Index + 1)
$namePara.Range.Text = "Prof. André Carvalhas"
